$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.486.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.93%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.512.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.51%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.577'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.538'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.66'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.112'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.899.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.545.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.97'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.675.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0970'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.94%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.39%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.97'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '250.01'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.76%  '

$ws.Range("E27").Value = '  +0.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.95'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '156.46'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0788'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.03%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.24'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.61%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.60%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.67%  '

$ws.Range("E38").Value = '  +1.92%  '

$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.25%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.119'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("E41").Value = '  -3.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0300'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.026.89'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.764.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.24%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.189'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '101.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.06%  '
